$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# A2: was a number (0), now an inline string "This is Delhi \n"
$ws.Range("A2").Value = "This is Delhi `n"
# The embedded newline makes Excel auto-grow the row height; re-fit it back
# down so we don't introduce an unrelated row-height change.
$ws.Rows(2).AutoFit()

# C2: was a number (21), now an inline string "55" (stored as text, not a
# number) - a leading apostrophe forces text entry, then reset the style so
# we don't leave the auto-applied "quote prefix" cell style behind.
$ws.Range("C2").Value = "'55"
$ws.Range("C2").Style = "Normal"
